# Revert "Change of position of rooms"
# Restructure the room map back to its earlier layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove cells that must disappear entirely (value + formatting)
# ---------------------------------------------------------------------
$ws.Range("C3").Clear() | Out-Null
$ws.Range("C7").Clear() | Out-Null
$ws.Range("A9").Clear() | Out-Null
$ws.Range("C9").Clear() | Out-Null
$ws.Range("A11").Clear() | Out-Null

# ---------------------------------------------------------------------
# 2) Plain value edits (style untouched - stays whatever it already is)
#    Order matters: the shared-string table records new unique strings
#    in first-write order, so "Starting_Room" and "Weight_Room" (which
#    must land right after "Dining_Room" in the table) are written
#    before any "Elevator" text, which is a brand-new string too but
#    must land at the very end of the table.
# ---------------------------------------------------------------------
$ws.Range("A2").Value2  = "Vault_Door"
$ws.Range("C2").Value2  = "Dining_Room"
$ws.Range("D2").Value2  = "Starting_Room"

$ws.Range("B3").Value2  = "."
$ws.Range("B4").Value2  = "."
$ws.Range("B5").Value2  = "."

$ws.Range("A6").Value2  = "Storage_Room"

$ws.Range("B7").Value2  = "."
$ws.Range("B8").Value2  = "."
$ws.Range("B9").Value2  = "."

$ws.Range("A10").Value2 = "Armory"
$ws.Range("C10").Value2 = "Hospital "

# ---------------------------------------------------------------------
# 3) Cells that need the small font (style index 1) applied from
#    scratch - they are brand-new cells in the extended range.
# ---------------------------------------------------------------------
$ws.Range("D10").Font.Size = 8
$ws.Range("D10").Value2 = "Science_Lab"

$ws.Range("B11").Font.Size = 8
$ws.Range("B11").Value2 = "."

$ws.Range("B12").Font.Size = 8
$ws.Range("B12").Value2 = "."

$ws.Range("B13").Font.Size = 8
$ws.Range("B13").Value2 = "."

$ws.Range("C14").Font.Size = 8
$ws.Range("C14").Value2 = "Weight_Room"

# "Elevator" is a brand-new string too, but must be the LAST new entry
# added to the shared-string table, so every occurrence is written here.
$ws.Range("B2").Value2  = "Elevator"
$ws.Range("B6").Value2  = "Elevator"
$ws.Range("B10").Value2 = "Elevator"

$ws.Range("B14").Font.Size = 8
$ws.Range("B14").Value2 = "Elevator"

# ---------------------------------------------------------------------
# 4) Cells that keep their text but become centered / get new text and
#    become centered. Clear C11's leftover text first.
# ---------------------------------------------------------------------
$ws.Range("C11").ClearContents() | Out-Null
$ws.Range("C5").Value2 = "Water_Treatment_Plant"
$ws.Range("C6").Value2 = "Power_Plant "

$ws.Range("C1:D1").HorizontalAlignment  = -4108
$ws.Range("C5:D5").HorizontalAlignment  = -4108
$ws.Range("C6:D6").HorizontalAlignment  = -4108
$ws.Range("C11:D11").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5) Merge the label cells (order matches the authored layout)
# ---------------------------------------------------------------------
$ws.Range("C1:D1").Merge()   | Out-Null
$ws.Range("C11:D11").Merge() | Out-Null
$ws.Range("C5:D5").Merge()   | Out-Null
$ws.Range("C6:D6").Merge()   | Out-Null

# ---------------------------------------------------------------------
# 6) Selection the author left active
# ---------------------------------------------------------------------
$ws.Range("B14").Select() | Out-Null
